$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 55.848606
$ws.Cells.Item(2, 8).Value = 167.545818
$ws.Cells.Item(2, 9).Value = 0.2323375192077237
$ws.Cells.Item(2, 10).Value = 0.2323375192077236
$ws.Cells.Item(2, 13).Value = 2.846571
$ws.Cells.Item(2, 14).Value = 8.539712999999999
$ws.Cells.Item(2, 15).Value = 0.01041928469143245
$ws.Cells.Item(2, 16).Value = 0.01041928469143244
$ws.Cells.Item(2, 17).Value = 158.977022230026
$ws.Cells.Item(2, 18).Value = 1430.793200070234
$ws.Cells.Item(2, 19).Value = 0.002420790757126427
$ws.Cells.Item(2, 20).Value = 0.002420790757126426
$ws.Cells.Item(3, 7).Value = 55.848606
$ws.Cells.Item(3, 8).Value = 167.545818
$ws.Cells.Item(3, 9).Value = 0.2323375192077237
$ws.Cells.Item(3, 10).Value = 0.2323375192077236
$ws.Cells.Item(3, 15).Value = 0.07835537840126532
$ws.Cells.Item(3, 16).Value = 0.0783553784012653
$ws.Cells.Item(3, 17).Value = 1195.543178139948
$ws.Cells.Item(3, 18).Value = 10759.88860325953
$ws.Cells.Item(3, 19).Value = 0.01820489423433243
$ws.Cells.Item(3, 20).Value = 0.01820489423433243
$ws.Cells.Item(4, 7).Value = 55.848606
$ws.Cells.Item(4, 8).Value = 167.545818
$ws.Cells.Item(4, 9).Value = 0.2323375192077237
$ws.Cells.Item(4, 10).Value = 0.2323375192077236
$ws.Cells.Item(4, 13).Value = 7.347547
$ws.Cells.Item(4, 14).Value = 22.042641
$ws.Cells.Item(4, 15).Value = 0.02689417688042223
$ws.Cells.Item(4, 16).Value = 0.02689417688042223
$ws.Cells.Item(4, 17).Value = 410.350257469482
$ws.Cells.Item(4, 18).Value = 3693.152317225338
$ws.Cells.Item(4, 19).Value = 0.006248526337531019
$ws.Cells.Item(4, 20).Value = 0.006248526337531017
$ws.Cells.Item(5, 7).Value = 55.848606
$ws.Cells.Item(5, 8).Value = 167.545818
$ws.Cells.Item(5, 9).Value = 0.2323375192077237
$ws.Cells.Item(5, 10).Value = 0.2323375192077236
$ws.Cells.Item(5, 13).Value = 241.601176
$ws.Cells.Item(5, 14).Value = 724.8035279999999
$ws.Cells.Item(5, 15).Value = 0.8843311600268801
$ws.Cells.Item(5, 16).Value = 0.8843311600268799
$ws.Cells.Item(5, 17).Value = 13493.08888756065
$ws.Cells.Item(5, 18).Value = 121437.7999880459
$ws.Cells.Item(5, 19).Value = 0.2054633078787338
$ws.Cells.Item(5, 20).Value = 0.2054633078787337
$ws.Cells.Item(6, 9).Value = 0.3515710112922583
$ws.Cells.Item(6, 10).Value = 0.3515710112922583
$ws.Cells.Item(6, 13).Value = 2.846571
$ws.Cells.Item(6, 14).Value = 8.539712999999999
$ws.Cells.Item(6, 15).Value = 0.01041928469143245
$ws.Cells.Item(6, 16).Value = 0.01041928469143244
$ws.Cells.Item(6, 17).Value = 240.5625775304569
$ws.Cells.Item(6, 18).Value = 2165.063197774113
$ws.Cells.Item(6, 19).Value = 0.00366311845590885
$ws.Cells.Item(6, 20).Value = 0.003663118455908849
$ws.Cells.Item(7, 9).Value = 0.3515710112922583
$ws.Cells.Item(7, 10).Value = 0.3515710112922583
$ws.Cells.Item(7, 15).Value = 0.07835537840126532
$ws.Cells.Item(7, 16).Value = 0.0783553784012653
$ws.Cells.Item(7, 19).Value = 0.02754747962472042
$ws.Cells.Item(7, 20).Value = 0.02754747962472041
$ws.Cells.Item(8, 9).Value = 0.3515710112922583
$ws.Cells.Item(8, 10).Value = 0.3515710112922583
$ws.Cells.Item(8, 13).Value = 7.347547
$ws.Cells.Item(8, 14).Value = 22.042641
$ws.Cells.Item(8, 15).Value = 0.02689417688042223
$ws.Cells.Item(8, 16).Value = 0.02689417688042223
$ws.Cells.Item(8, 17).Value = 620.9382604003822
$ws.Cells.Item(8, 18).Value = 5588.444343603441
$ws.Cells.Item(8, 19).Value = 0.009455212963722917
$ws.Cells.Item(8, 20).Value = 0.009455212963722915
$ws.Cells.Item(9, 9).Value = 0.3515710112922583
$ws.Cells.Item(9, 10).Value = 0.3515710112922583
$ws.Cells.Item(9, 13).Value = 241.601176
$ws.Cells.Item(9, 14).Value = 724.8035279999999
$ws.Cells.Item(9, 15).Value = 0.8843311600268801
$ws.Cells.Item(9, 16).Value = 0.8843311600268799
$ws.Cells.Item(9, 17).Value = 20417.61882382332
$ws.Cells.Item(9, 18).Value = 183758.5694144099
$ws.Cells.Item(9, 19).Value = 0.3109052002479061
$ws.Cells.Item(9, 20).Value = 0.3109052002479061
$ws.Cells.Item(10, 7).Value = 33.195992
$ws.Cells.Item(10, 8).Value = 99.58797600000001
$ws.Cells.Item(10, 9).Value = 0.1380996766314891
$ws.Cells.Item(10, 10).Value = 0.1380996766314891
$ws.Cells.Item(10, 13).Value = 2.846571
$ws.Cells.Item(10, 14).Value = 8.539712999999999
$ws.Cells.Item(10, 15).Value = 0.01041928469143245
$ws.Cells.Item(10, 16).Value = 0.01041928469143244
$ws.Cells.Item(10, 17).Value = 94.49474814343199
$ws.Cells.Item(10, 18).Value = 850.452733290888
$ws.Cells.Item(10, 19).Value = 0.001438899846618246
$ws.Cells.Item(10, 20).Value = 0.001438899846618245
$ws.Cells.Item(11, 7).Value = 33.195992
$ws.Cells.Item(11, 8).Value = 99.58797600000001
$ws.Cells.Item(11, 9).Value = 0.1380996766314891
$ws.Cells.Item(11, 10).Value = 0.1380996766314891
$ws.Cells.Item(11, 15).Value = 0.07835537840126532
$ws.Cells.Item(11, 16).Value = 0.0783553784012653
$ws.Cells.Item(11, 17).Value = 710.6218869131361
$ws.Cells.Item(11, 18).Value = 6395.596982218225
$ws.Cells.Item(11, 19).Value = 0.01082085241955271
$ws.Cells.Item(11, 20).Value = 0.0108208524195527
$ws.Cells.Item(12, 7).Value = 33.195992
$ws.Cells.Item(12, 8).Value = 99.58797600000001
$ws.Cells.Item(12, 9).Value = 0.1380996766314891
$ws.Cells.Item(12, 10).Value = 0.1380996766314891
$ws.Cells.Item(12, 13).Value = 7.347547
$ws.Cells.Item(12, 14).Value = 22.042641
$ws.Cells.Item(12, 15).Value = 0.02689417688042223
$ws.Cells.Item(12, 16).Value = 0.02689417688042223
$ws.Cells.Item(12, 17).Value = 243.909111431624
$ws.Cells.Item(12, 18).Value = 2195.182002884616
$ws.Cells.Item(12, 19).Value = 0.003714077130456381
$ws.Cells.Item(12, 20).Value = 0.00371407713045638
$ws.Cells.Item(13, 7).Value = 33.195992
$ws.Cells.Item(13, 8).Value = 99.58797600000001
$ws.Cells.Item(13, 9).Value = 0.1380996766314891
$ws.Cells.Item(13, 10).Value = 0.1380996766314891
$ws.Cells.Item(13, 13).Value = 241.601176
$ws.Cells.Item(13, 14).Value = 724.8035279999999
$ws.Cells.Item(13, 15).Value = 0.8843311600268801
$ws.Cells.Item(13, 16).Value = 0.8843311600268799
$ws.Cells.Item(13, 17).Value = 8020.190705686592
$ws.Cells.Item(13, 18).Value = 72181.71635117933
$ws.Cells.Item(13, 19).Value = 0.1221258472348618
$ws.Cells.Item(13, 20).Value = 0.1221258472348618
$ws.Cells.Item(14, 7).Value = 66.82284533333335
$ws.Cells.Item(14, 8).Value = 200.468536
$ws.Cells.Item(14, 9).Value = 0.277991792868529
$ws.Cells.Item(14, 10).Value = 0.2779917928685289
$ws.Cells.Item(14, 13).Value = 2.846571
$ws.Cells.Item(14, 14).Value = 8.539712999999999
$ws.Cells.Item(14, 15).Value = 0.01041928469143245
$ws.Cells.Item(14, 16).Value = 0.01041928469143244
$ws.Cells.Item(14, 17).Value = 190.215973663352
$ws.Cells.Item(14, 18).Value = 1711.943762970168
$ws.Cells.Item(14, 19).Value = 0.002896475631778923
$ws.Cells.Item(14, 20).Value = 0.002896475631778922
$ws.Cells.Item(15, 7).Value = 66.82284533333335
$ws.Cells.Item(15, 8).Value = 200.468536
$ws.Cells.Item(15, 9).Value = 0.277991792868529
$ws.Cells.Item(15, 10).Value = 0.2779917928685289
$ws.Cells.Item(15, 15).Value = 0.07835537840126532
$ws.Cells.Item(15, 16).Value = 0.0783553784012653
$ws.Cells.Item(15, 17).Value = 1430.46716120663
$ws.Cells.Item(15, 18).Value = 12874.20445085967
$ws.Cells.Item(15, 19).Value = 0.02178215212265976
$ws.Cells.Item(15, 20).Value = 0.02178215212265975
$ws.Cells.Item(16, 7).Value = 66.82284533333335
$ws.Cells.Item(16, 8).Value = 200.468536
$ws.Cells.Item(16, 9).Value = 0.277991792868529
$ws.Cells.Item(16, 10).Value = 0.2779917928685289
$ws.Cells.Item(16, 13).Value = 7.347547
$ws.Cells.Item(16, 14).Value = 22.042641
$ws.Cells.Item(16, 15).Value = 0.02689417688042223
$ws.Cells.Item(16, 16).Value = 0.02689417688042223
$ws.Cells.Item(16, 17).Value = 490.9839967603974
$ws.Cells.Item(16, 18).Value = 4418.855970843577
$ws.Cells.Item(16, 19).Value = 0.007476360448711919
$ws.Cells.Item(16, 20).Value = 0.007476360448711916
$ws.Cells.Item(17, 7).Value = 66.82284533333335
$ws.Cells.Item(17, 8).Value = 200.468536
$ws.Cells.Item(17, 9).Value = 0.277991792868529
$ws.Cells.Item(17, 10).Value = 0.2779917928685289
$ws.Cells.Item(17, 13).Value = 241.601176
$ws.Cells.Item(17, 14).Value = 724.8035279999999
$ws.Cells.Item(17, 15).Value = 0.8843311600268801
$ws.Cells.Item(17, 16).Value = 0.8843311600268799
$ws.Cells.Item(17, 17).Value = 16144.47801619945
$ws.Cells.Item(17, 18).Value = 145300.302145795
$ws.Cells.Item(17, 19).Value = 0.2458368046653784
$ws.Cells.Item(17, 20).Value = 0.2458368046653783

Write-Output "done"